$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("D:D").Insert()
Write-Host $ws.Range("D7").Value2
Write-Host $ws.Range("E7").Value2
Write-Host $ws.Range("K7").Value2
Write-Host $ws.Range("L7").Value2
